$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-7 (old extra schedule rows) so only rows 1-3 remain
$ws.Rows("4:7").Delete()

# Row 1: name / role / year
$ws.Range("A1").Value = "ทดสอบ"
$ws.Range("B1").Value = "คุณครู"
$ws.Range("C1").Value = 2561

# Row 2: headers (unchanged content)
$ws.Range("A2").Value = "date"
$ws.Range("B2").Value = "เวลาในการสอน"
$ws.Range("C2").Value = "ชั้น"
$ws.Range("D2").Value = "รหัสวิชา"
$ws.Range("E2").Value = "ระดับชั้นเรียน"

# Row 3: new single data row
$ws.Range("A3").Value = "วันพุธ"
$ws.Range("B3").Value = "14:00 - 15:00"
$ws.Range("C3").Value = "ป.1/1"
$ws.Range("D3").Value = "ค 11101"
$ws.Range("E3").Value = "ประถมศึกษา"

$ws.Range("K12").Select() | Out-Null
